$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# "Valor Mora" (overdue amount) for period 2106 (row 16) and period 2011 (row 23)
# are swapped as part of the account-statement database update.
$ws.Range("F16").Value = 35112
$ws.Range("F23").Value = 24578
